$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Title / absolute-path metadata updates (October -> November 2016 refresh)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Table 1.1. Net Generation by Energy Source:  Total (All Sectors), 2006-November 2016"

# ---------------------------------------------------------------------------
# 2. Shift the "Year to Date" / "Rolling 12 Months" block down by one row to
#    make room for a new "November" monthly-data row at row 53.
#    Work bottom-up so the source of each copy is read before being
#    overwritten.
# ---------------------------------------------------------------------------
$ws.Range("A60:P60").Copy($ws.Range("A61:P61"))
$ws.Rows.Item(61).RowHeight = 237.95

$ws.Range("A59:P59").Copy($ws.Range("A60:P60"))
$ws.Rows.Item(60).AutoFit()

$ws.Range("A58:P58").Copy($ws.Range("A59:P59"))
$ws.Range("A57:P57").Copy($ws.Range("A58:P58"))
$ws.Range("A56:P56").Copy($ws.Range("A57:P57"))
$ws.Range("A55:P55").Copy($ws.Range("A56:P56"))
$ws.Range("A54:P54").Copy($ws.Range("A55:P55"))
$ws.Range("A53:P53").Copy($ws.Range("A54:P54"))

# ---------------------------------------------------------------------------
# 3. New row 53: November monthly data (style copied from the October row so
#    the existing number/label formats -- not new ones -- are reused).
# ---------------------------------------------------------------------------
$ws.Range("A52:P52").Copy($ws.Range("A53:P53"))
$ws.Range("A53").UnMerge()

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 87000
$ws.Range("C53").Value = 1058
$ws.Range("D53").Value = 781
$ws.Range("E53").Value = 94586
$ws.Range("F53").Value = 1001
$ws.Range("G53").Value = 65179
$ws.Range("H53").Value = 18815
$ws.Range("I53").Value = 2642
$ws.Range("J53").Value = 25874
$ws.Range("K53").Value = -607
$ws.Range("L53").Value = 1093
$ws.Range("M53").Value = 297422
$ws.Range("N53").Value = 1307
$ws.Range("O53").Value = 3766
$ws.Range("P53").Value = 3950

# ---------------------------------------------------------------------------
# 4. "Year to Date" figures now cover Jan-Nov instead of Jan-Oct.
# ---------------------------------------------------------------------------
# Row 55 = 2014
$ws.Range("B55").Value = 1457090
$ws.Range("C55").Value = 17330
$ws.Range("D55").Value = 10807
$ws.Range("E55").Value = 1035571
$ws.Range("F55").Value = 10869
$ws.Range("G55").Value = 723803
$ws.Range("H55").Value = 237038
$ws.Range("I55").Value = 16659
$ws.Range("J55").Value = 239932
$ws.Range("K55").Value = -5694
$ws.Range("L55").Value = 12244
$ws.Range("M55").Value = 3755649
$ws.Range("N55").Value = 10467
$ws.Range("O55").Value = 24780
$ws.Range("P55").Value = 27126

# Row 56 = 2015
$ws.Range("B56").Value = 1262903
$ws.Range("C56").Value = 16424
$ws.Range("D56").Value = 10128
$ws.Range("E56").Value = 1223705
$ws.Range("F56").Value = 12007
$ws.Range("G56").Value = 727544
$ws.Range("H56").Value = 225915
$ws.Range("I56").Value = 23323
$ws.Range("J56").Value = 243236
$ws.Range("K56").Value = -4811
$ws.Range("L56").Value = 12800
$ws.Range("M56").Value = 3753174
$ws.Range("N56").Value = 13225
$ws.Range("O56").Value = 33447
$ws.Range("P56").Value = 36548

# Row 57 = 2016
$ws.Range("B57").Value = 1121120
$ws.Range("C57").Value = 11528
$ws.Range("D57").Value = 10366
$ws.Range("E57").Value = 1284457
$ws.Range("F57").Value = 11987
$ws.Range("G57").Value = 733632
$ws.Range("H57").Value = 243220
$ws.Range("I57").Value = 33832
$ws.Range("J57").Value = 276232
$ws.Range("K57").Value = -5933
$ws.Range("L57").Value = 12550
$ws.Range("M57").Value = 3732992
$ws.Range("N57").Value = 18281
$ws.Range("O57").Value = 48820
$ws.Range("P57").Value = 52113

# ---------------------------------------------------------------------------
# 5. "Rolling 12 Months Ending in October" -> "...Ending in November", and
#    its two data rows refresh to the Dec 2015-Nov 2016 / Dec 2014-Nov 2015
#    windows respectively.
# ---------------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Row 59 = 2015
$ws.Range("B59").Value = 1387523
$ws.Range("C59").Value = 17371
$ws.Range("D59").Value = 11276
$ws.Range("E59").Value = 1314743
$ws.Range("F59").Value = 13160
$ws.Range("G59").Value = 800907
$ws.Range("H59").Value = 248243
$ws.Range("I59").Value = 24355
$ws.Range("J59").Value = 264826
$ws.Range("K59").Value = -5291
$ws.Range("L59").Value = 14017
$ws.Range("M59").Value = 4091130
$ws.Range("N59").Value = 13991
$ws.Range("O59").Value = 35150
$ws.Range("P59").Value = 38346

# Row 60 = 2016
$ws.Range("B60").Value = 1210615
$ws.Range("C60").Value = 12475
$ws.Range("D60").Value = 11115
$ws.Range("E60").Value = 1394234
$ws.Range("F60").Value = 13097
$ws.Range("G60").Value = 803266
$ws.Range("H60").Value = 266386
$ws.Range("I60").Value = 35402
$ws.Range("J60").Value = 303264
$ws.Range("K60").Value = -6214
$ws.Range("L60").Value = 13779
$ws.Range("M60").Value = 4057419
$ws.Range("N60").Value = 19195
$ws.Range("O60").Value = 51179
$ws.Range("P60").Value = 54597
